$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 256
$ws1.Range("F6").Value = 593
$ws1.Range("F9").Value = 281
$ws1.Range("F10").Value = 406
$ws1.Range("F12").Value = 1165
$ws1.Range("F15").Value = 13
$ws1.Range("F16").Value = 1553
$ws1.Range("F17").Value = 1553
$ws1.Range("F18").Value = 1282
$ws1.Range("F20").Value = 1369
$ws1.Range("F22").Value = 372
$ws1.Range("F26").Value = 6776
$ws1.Range("F27").Value = 5547
$ws1.Range("F29").Value = 157
$ws1.Range("F30").Value = 490
$ws1.Range("F40").Value = 637
$ws1.Range("F41").Value = 20
$ws1.Range("F43").Value = 278
$ws1.Range("F45").Value = 161
$ws1.Range("F46").Value = 68
$ws1.Range("G46").Value = 76
$ws1.Range("F47").Value = 94
$ws1.Range("F48").Value = 108
$ws1.Range("F49").Value = 10

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 256

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2487
$ws3.Range("F4").Value = 224

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 256
$ws4.Range("F6").Value = 224
$ws4.Range("F9").Value = 593
$ws4.Range("F12").Value = 281
$ws4.Range("F14").Value = 406
$ws4.Range("F16").Value = 1168
$ws4.Range("F19").Value = 13
$ws4.Range("F20").Value = 1553
$ws4.Range("F21").Value = 1553
$ws4.Range("F22").Value = 1282
$ws4.Range("F24").Value = 372
$ws4.Range("F29").Value = 6776
$ws4.Range("F30").Value = 5551
$ws4.Range("F38").Value = 637
$ws4.Range("F43").Value = 278
$ws4.Range("F44").Value = 161
$ws4.Range("F45").Value = 68
$ws4.Range("G45").Value = 76
$ws4.Range("F46").Value = 94
$ws4.Range("F47").Value = 108
$ws4.Range("F49").Value = 256
$ws4.Range("F50").Value = 10
